$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains its text representation (values such as
# "28.944.54" or "0.9982" must not be auto-converted to numbers/dates).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.944.54"
$ws.Range("E2").Value = "  +5.38%  "

$ws.Range("D3").Value = "1.913.89"
$ws.Range("E3").Value = "  +4.84%  "

$ws.Range("D4").Value = "0.9982"
$ws.Range("E4").Value = "  -0.75%  "

$ws.Range("D5").Value = "338.74"
$ws.Range("E5").Value = "  +1.83%  "

$ws.Range("D6").Value = "0.9986"
$ws.Range("E6").Value = "  -0.63%  "

$ws.Range("D7").Value = "0.4734"
$ws.Range("E7").Value = "  +3.41%  "

$ws.Range("D8").Value = "0.4057"
$ws.Range("E8").Value = "  +6.70%  "

$ws.Range("D9").Value = "47.97"
$ws.Range("E9").Value = "  +3.47%  "

$ws.Range("D10").Value = "0.08177"
$ws.Range("E10").Value = "  +3.82%  "

$ws.Range("D11").Value = "1.025"
$ws.Range("E11").Value = "  +5.78%  "

$ws.Range("D12").Value = "22.54"
$ws.Range("E12").Value = "  +7.30%  "

$ws.Range("D13").Value = "1.902.71"
$ws.Range("E13").Value = "  +4.16%  "

$ws.Range("D14").Value = "6.095"

$ws.Range("D15").Value = "7.367"
$ws.Range("E15").Value = "  +4.59%  "

$ws.Range("D16").Value = "91.66"
$ws.Range("E16").Value = "  +2.34%  "

$ws.Range("D17").Value = "0.9986"
$ws.Range("E17").Value = "  -0.73%  "

$ws.Range("D18").Value = "0.00001054"
$ws.Range("E18").Value = "  +2.82%  "

$ws.Range("D19").Value = "0.06623"
$ws.Range("E19").Value = "  -0.31%  "

$ws.Range("D20").Value = "17.88"
$ws.Range("E20").Value = "  +4.67%  "

$ws.Range("D21").Value = "0.9983"
$ws.Range("E21").Value = "  -0.71%  "

$ws.Range("D22").Value = "28.954.67"
$ws.Range("E22").Value = "  +5.50%  "

$ws.Range("D23").Value = "5.569"
$ws.Range("E23").Value = "  +4.47%  "

$ws.Range("D24").Value = "11.18"
$ws.Range("E24").Value = "  +3.44%  "

$ws.Range("D25").Value = "2.266"
$ws.Range("E25").Value = "  -1.84%  "

$ws.Range("D26").Value = "2.118.48"
$ws.Range("E26").Value = "  +3.71%  "

$ws.Range("D27").Value = "161.18"
$ws.Range("E27").Value = "  +3.53%  "

$ws.Range("D28").Value = "20.08"
$ws.Range("E28").Value = "  +3.67%  "

$ws.Range("D29").Value = "2.185"
$ws.Range("E29").Value = "  +6.47%  "

$ws.Range("D30").Value = "5.539"
$ws.Range("E30").Value = "  +5.11%  "

$ws.Range("D31").Value = "121.03"
$ws.Range("E31").Value = "  +2.38%  "

$ws.Range("D32").Value = "1.020"
$ws.Range("E32").Value = "  +8.31%  "

$ws.Range("D33").Value = "0.09587"
$ws.Range("E33").Value = "  +3.07%  "

$ws.Range("D34").Value = "3.650"
$ws.Range("E34").Value = "  +1.68%  "

$ws.Range("E35").Value = "  +7.33%  "

$ws.Range("D36").Value = "5.432"
$ws.Range("E36").Value = "  +3.52%  "

$ws.Range("D37").Value = "0.06214"
$ws.Range("E37").Value = "  +5.04%  "

$ws.Range("D38").Value = "0.02291"
$ws.Range("E38").Value = "  +5.25%  "

$ws.Range("D39").Value = "8.668"
$ws.Range("E39").Value = "  +7.69%  "

$ws.Range("D40").Value = "1.200"
$ws.Range("E40").Value = "  +3.77%  "

$ws.Range("D41").Value = "0.6052"
$ws.Range("E41").Value = "  +5.14%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "10.63"
$ws.Range("E42").Value = "  +6.56%  "

$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").Value = "0.1906"
$ws.Range("E43").Value = "  +4.39%  "

$ws.Range("D44").Value = "0.9985"
$ws.Range("E44").Value = "  -0.59%  "

$ws.Range("D45").Value = "1.285"
$ws.Range("E45").Value = "  +1.40%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.5655"
$ws.Range("E46").Value = "  +4.02%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "12.38"
$ws.Range("E47").Value = "  +3.89%  "

$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "1.991"
$ws.Range("E48").Value = "  +6.85%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "0.07296"
$ws.Range("E49").Value = "  +10.50%  "

$ws.Range("D50").Value = "2.145"
$ws.Range("E50").Value = "  +17.95%  "

$ws.Range("D51").Value = "113.26"
$ws.Range("E51").Value = "  +2.22%  "
